$d = $word.ActiveDocument

# Update the date/weekday heading paragraph.
$d.Content.Find.Execute("2024-02-28 Wednesday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-02-29 Thursday", 2)

# Update the division-fact table cells (row, col) -> new text.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "25÷8=3, 1" },
    @{ Row = 1;  Col = 2; Text = "67÷7=9, 4" },
    @{ Row = 1;  Col = 3; Text = "44÷5=8, 4" },
    @{ Row = 1;  Col = 4; Text = "11÷2=5, 1" },
    @{ Row = 1;  Col = 5; Text = "41÷7=5, 6" },

    @{ Row = 5;  Col = 1; Text = "78÷3=26, 0" },
    @{ Row = 5;  Col = 2; Text = "69÷7=9, 6" },
    @{ Row = 5;  Col = 3; Text = "28÷9=3, 1" },
    @{ Row = 5;  Col = 4; Text = "60÷6=10, 0" },
    @{ Row = 5;  Col = 5; Text = "27÷5=5, 2" },

    @{ Row = 9;  Col = 1; Text = "97÷3=32, 1" },
    @{ Row = 9;  Col = 2; Text = "39÷9=4, 3" },
    @{ Row = 9;  Col = 3; Text = "86÷4=21, 2" },
    @{ Row = 9;  Col = 4; Text = "34÷2=17, 0" },
    @{ Row = 9;  Col = 5; Text = "36÷9=4, 0" },

    @{ Row = 13; Col = 1; Text = "58÷6=9, 4" },
    @{ Row = 13; Col = 2; Text = "55÷8=6, 7" },
    @{ Row = 13; Col = 3; Text = "40÷5=8, 0" },
    @{ Row = 13; Col = 4; Text = "27÷7=3, 6" },
    @{ Row = 13; Col = 5; Text = "21÷9=2, 3" },

    @{ Row = 17; Col = 1; Text = "77÷4=19, 1" },
    @{ Row = 17; Col = 2; Text = "81÷4=20, 1" },
    @{ Row = 17; Col = 3; Text = "85÷3=28, 1" },
    @{ Row = 17; Col = 4; Text = "91÷7=13, 0" },
    @{ Row = 17; Col = 5; Text = "44÷7=6, 2" }
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
